$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 16
$ws.Range("F4").Value = 187
$ws.Range("F5").Value = 1028
$ws.Range("F7").Value = 2626
$ws.Range("F9").Value = 1280
$ws.Range("F10").Value = 922
$ws.Range("F11").Value = 618
$ws.Range("F12").Value = 931
$ws.Range("F13").Value = 1168
$ws.Range("F17").Value = 786
$ws.Range("F18").Value = 223
$ws.Range("F19").Value = 517
$ws.Range("F20").Value = 1130
$ws.Range("F23").Value = 604
$ws.Range("F24").Value = 223
$ws.Range("F25").Value = 311
$ws.Range("F27").Value = 691
$ws.Range("F28").Value = 550
$ws.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202407/Hllpq7tZ1722419749368.jpeg"
$ws.Range("F29").Value = 5604
$ws.Range("F30").Value = 5604
$ws.Range("F35").Value = 170
$ws.Range("F36").Value = 1630
$ws.Range("F37").Value = 9
$ws.Range("F38").Value = 84
$ws.Range("F41").Value = 89
$ws.Range("F42").Value = 148
$ws.Range("F43").Value = 5
$ws.Range("F46").Value = 140
$ws.Range("F47").Value = 138
$ws.Range("F48").Value = 114
$ws.Range("F49").Value = 25
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 85
$ws.Range("F5").Value = 8
$ws.Range("F9").Value = 36
$ws.Range("F10").Value = 180
$ws.Range("F12").Value = 194
$ws.Range("F13").Value = 4410
$ws.Range("F14").Value = 33
$ws.Range("F17").Value = 35
$ws.Range("F18").Value = 206
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2279
$ws.Range("F3").Value = 737
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2279
$ws.Range("F3").Value = 16
$ws.Range("F4").Value = 85
$ws.Range("F5").Value = 1028
$ws.Range("F6").Value = 2626
$ws.Range("F8").Value = 1280
$ws.Range("F9").Value = 922
$ws.Range("F10").Value = 618
$ws.Range("F11").Value = 931
$ws.Range("F12").Value = 1168
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 786
$ws.Range("F19").Value = 223
$ws.Range("F20").Value = 517
$ws.Range("F21").Value = 1130
$ws.Range("F25").Value = 36
$ws.Range("F27").Value = 604
$ws.Range("F28").Value = 311
$ws.Range("F30").Value = 550
$ws.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202407/Hllpq7tZ1722419749368.jpeg"
$ws.Range("F31").Value = 5604
$ws.Range("F32").Value = 194
$ws.Range("F36").Value = 170
$ws.Range("F37").Value = 1630
$ws.Range("F38").Value = 9
$ws.Range("F40").Value = 33
$ws.Range("F41").Value = 33
$ws.Range("F43").Value = 89
$ws.Range("F44").Value = 35
$ws.Range("F45").Value = 148
$ws.Range("F49").Value = 114
